$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.936.48"
$ws.Range("E2").Value = "  -3.51%  "
$ws.Range("D3").Value = "'2.918.72"
$ws.Range("E3").Value = "  -4.13%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'586.02"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'145.47"
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").Value = "'2.916.69"
$ws.Range("E9").Value = "  -4.22%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -5.13%  "
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "'33.60"
$ws.Range("E14").Value = "  -6.17%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'3.402.75"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").Value = "'60.875.46"
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Value = "'6.77"
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("D19").Value = "'2.921.74"
$ws.Range("E19").Value = "  -3.96%  "
$ws.Range("D20").Value = "'429.79"
$ws.Range("E20").Value = "  -5.79%  "
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("D22").Value = "'0.681"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  -5.56%  "
$ws.Range("D24").Value = "'80.46"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("E33").Value = "  -4.27%  "
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "'0.0₃0868"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("E37").Value = "  -5.24%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.01"
$ws.Range("E38").Value = "  -6.45%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.127"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").Value = "'49.47"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "'2.00"
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("D42").Value = "'8.63"
$ws.Range("E42").Value = "  -5.56%  "
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("D44").Value = "'41.50"
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").Value = "'377.68"
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").Value = "'2.698.97"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'132.77"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D50").Value = "'24.88"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  -2.42%  "
